$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; insert it as a new row 66
# (pushing the existing rows 66-97 down to 67-98) and fill it in with the
# new reading. All the "descriptive" columns (market, region, product,
# unit, origin, classification, etc.) repeat the same values used by every
# other row for this product/market.
$ws.Rows("66:66").Insert()

$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44795
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = 100112038
$ws.Range("G66").Value = "Cebollín baby"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 300
$ws.Range("K66").Value = 2000
$ws.Range("L66").Value = 2500
$ws.Range("M66").Value = 2250
$ws.Range("N66").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O66").Value = "Región de Arica y Parinacota"
$ws.Range("P66").Value = 1125
$ws.Range("Q66").Value = 2
$ws.Range("R66").Value = "Hortaliza"
